# The workbook has one sheet per property type. The "property_category"
# column on each sheet was mistakenly left as "land" (copied from the
# 土地/Land sheet) when the sheet was generated. This fixes the
# property_category values for the 建物 (Building) and 汽車 (Car) sheets
# so they correctly read "building" and "car" respectively.

$wb = $excel.ActiveWorkbook

# 建物 (Building) sheet: column I is "property_category" for data rows 2-15.
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2:I15").Value = "building"

# 汽車 (Car) sheet: column H is "property_category" for data row 2.
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
